$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 68, shifting existing rows 68:169 down to 69:170.
$ws.Rows(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 44771
$ws.Range("D68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = 100112026
$ws.Range("G68").Value = "Haba"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 100
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 16000
$ws.Range("M68").Value = 15450
$ws.Range("N68").Value = "$/saco 25 kilos"
$ws.Range("O68").Value = "Provincia de Limarí"
$ws.Range("P68").Value = 618
$ws.Range("Q68").Value = 25
$ws.Range("R68").Value = "Hortaliza"
